$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value = 0.4549446666666667
$ws.Cells.Item(2, 8).Value = 1.364834
$ws.Cells.Item(2, 9).Value = 0.8656500014587819
$ws.Cells.Item(2, 10).Value = 0.8656500014587818
$ws.Cells.Item(2, 13).Value = 9.101967
$ws.Cells.Item(2, 14).Value = 27.305901
$ws.Cells.Item(2, 15).Value = 0.2872601673725235
$ws.Cells.Item(2, 16).Value = 0.2872601673725235
$ws.Cells.Item(2, 17).Value = 4.140891342826
$ws.Cells.Item(2, 18).Value = 37.268022085434
$ws.Cells.Item(2, 19).Value = 0.2486667643050749
$ws.Cells.Item(2, 20).Value = 0.2486667643050748

# Row 3
$ws.Cells.Item(3, 7).Value = 0.4549446666666667
$ws.Cells.Item(3, 8).Value = 1.364834
$ws.Cells.Item(3, 9).Value = 0.8656500014587819
$ws.Cells.Item(3, 10).Value = 0.8656500014587818
$ws.Cells.Item(3, 15).Value = 0.3055950511371977
$ws.Cells.Item(3, 16).Value = 0.3055950511371977
$ws.Cells.Item(3, 17).Value = 4.405190991981334
$ws.Cells.Item(3, 18).Value = 39.646718927832
$ws.Cells.Item(3, 19).Value = 0.2645383564627117
$ws.Cells.Item(3, 20).Value = 0.2645383564627117

# Row 4
$ws.Cells.Item(4, 7).Value = 0.4549446666666667
$ws.Cells.Item(4, 8).Value = 1.364834
$ws.Cells.Item(4, 9).Value = 0.8656500014587819
$ws.Cells.Item(4, 10).Value = 0.8656500014587818
$ws.Cells.Item(4, 13).Value = 3.905093666666666
$ws.Cells.Item(4, 14).Value = 11.715281
$ws.Cells.Item(4, 15).Value = 0.1232456523180152
$ws.Cells.Item(4, 16).Value = 0.1232456523180152
$ws.Cells.Item(4, 17).Value = 1.776601536483778
$ws.Cells.Item(4, 18).Value = 15.989413828354
$ws.Cells.Item(4, 19).Value = 0.1066875991088784
$ws.Cells.Item(4, 20).Value = 0.1066875991088784

# Row 5
$ws.Cells.Item(5, 7).Value = 0.4549446666666667
$ws.Cells.Item(5, 8).Value = 1.364834
$ws.Cells.Item(5, 9).Value = 0.8656500014587819
$ws.Cells.Item(5, 10).Value = 0.8656500014587818
$ws.Cells.Item(5, 13).Value = 6.285238333333333
$ws.Cells.Item(5, 14).Value = 18.855715
$ws.Cells.Item(5, 15).Value = 0.1983635642284282
$ws.Cells.Item(5, 16).Value = 0.1983635642284282
$ws.Cells.Item(5, 17).Value = 2.859435658478889
$ws.Cells.Item(5, 18).Value = 25.73492092631
$ws.Cells.Item(5, 19).Value = 0.171713419663708
$ws.Cells.Item(5, 20).Value = 0.171713419663708

# Row 6
$ws.Cells.Item(6, 7).Value = 0.4549446666666667
$ws.Cells.Item(6, 8).Value = 1.364834
$ws.Cells.Item(6, 9).Value = 0.8656500014587819
$ws.Cells.Item(6, 10).Value = 0.8656500014587818
$ws.Cells.Item(6, 13).Value = 2.710232666666667
$ws.Cells.Item(6, 14).Value = 8.130698000000001
$ws.Cells.Item(6, 15).Value = 0.08553556494383548
$ws.Cells.Item(6, 16).Value = 0.08553556494383548
$ws.Cells.Item(6, 17).Value = 1.233005897125778
$ws.Cells.Item(6, 18).Value = 11.097053074132
$ws.Cells.Item(6, 19).Value = 0.07404386191840892
$ws.Cells.Item(6, 20).Value = 0.07404386191840891

# Row 7
$ws.Cells.Item(7, 9).Value = 0.02407814503842938
$ws.Cells.Item(7, 10).Value = 0.02407814503842938
$ws.Cells.Item(7, 13).Value = 9.101967
$ws.Cells.Item(7, 14).Value = 27.305901
$ws.Cells.Item(7, 15).Value = 0.2872601673725235
$ws.Cells.Item(7, 16).Value = 0.2872601673725235
$ws.Cells.Item(7, 17).Value = 0.115179324407
$ws.Cells.Item(7, 18).Value = 1.036613919663
$ws.Cells.Item(7, 19).Value = 0.00691669197375912
$ws.Cells.Item(7, 20).Value = 0.00691669197375912

# Row 8
$ws.Cells.Item(8, 9).Value = 0.02407814503842938
$ws.Cells.Item(8, 10).Value = 0.02407814503842938
$ws.Cells.Item(8, 15).Value = 0.3055950511371977
$ws.Cells.Item(8, 16).Value = 0.3055950511371977
$ws.Cells.Item(8, 19).Value = 0.00735816196430769
$ws.Cells.Item(8, 20).Value = 0.00735816196430769

# Row 9
$ws.Cells.Item(9, 9).Value = 0.02407814503842938
$ws.Cells.Item(9, 10).Value = 0.02407814503842938
$ws.Cells.Item(9, 13).Value = 3.905093666666666
$ws.Cells.Item(9, 14).Value = 11.715281
$ws.Cells.Item(9, 15).Value = 0.1232456523180152
$ws.Cells.Item(9, 16).Value = 0.1232456523180152
$ws.Cells.Item(9, 17).Value = 0.04941635695588888
$ws.Cells.Item(9, 18).Value = 0.4447472126029999
$ws.Cells.Item(9, 19).Value = 0.002967526691869011
$ws.Cells.Item(9, 20).Value = 0.002967526691869011

# Row 10
$ws.Cells.Item(10, 9).Value = 0.02407814503842938
$ws.Cells.Item(10, 10).Value = 0.02407814503842938
$ws.Cells.Item(10, 13).Value = 6.285238333333333
$ws.Cells.Item(10, 14).Value = 18.855715
$ws.Cells.Item(10, 15).Value = 0.1983635642284282
$ws.Cells.Item(10, 16).Value = 0.1983635642284282
$ws.Cells.Item(10, 17).Value = 0.07953550094944443
$ws.Cells.Item(10, 18).Value = 0.715819508545
$ws.Cells.Item(10, 19).Value = 0.004776226669831896
$ws.Cells.Item(10, 20).Value = 0.004776226669831896

# Row 11
$ws.Cells.Item(11, 9).Value = 0.02407814503842938
$ws.Cells.Item(11, 10).Value = 0.02407814503842938
$ws.Cells.Item(11, 13).Value = 2.710232666666667
$ws.Cells.Item(11, 14).Value = 8.130698000000001
$ws.Cells.Item(11, 15).Value = 0.08553556494383548
$ws.Cells.Item(11, 16).Value = 0.08553556494383548
$ws.Cells.Item(11, 17).Value = 0.03429618757488889
$ws.Cells.Item(11, 18).Value = 0.308665688174
$ws.Cells.Item(11, 19).Value = 0.002059537738661666
$ws.Cells.Item(11, 20).Value = 0.002059537738661666

# Row 12
$ws.Cells.Item(12, 7).Value = 0.05795366666666666
$ws.Cells.Item(12, 8).Value = 0.173861
$ws.Cells.Item(12, 9).Value = 0.1102718535027888
$ws.Cells.Item(12, 10).Value = 0.1102718535027888
$ws.Cells.Item(12, 13).Value = 9.101967
$ws.Cells.Item(12, 14).Value = 27.305901
$ws.Cells.Item(12, 15).Value = 0.2872601673725235
$ws.Cells.Item(12, 16).Value = 0.2872601673725235
$ws.Cells.Item(12, 17).Value = 0.5274923615289999
$ws.Cells.Item(12, 18).Value = 4.747431253761
$ws.Cells.Item(12, 19).Value = 0.0316767110936895
$ws.Cells.Item(12, 20).Value = 0.0316767110936895

# Row 13
$ws.Cells.Item(13, 7).Value = 0.05795366666666666
$ws.Cells.Item(13, 8).Value = 0.173861
$ws.Cells.Item(13, 9).Value = 0.1102718535027888
$ws.Cells.Item(13, 10).Value = 0.1102718535027888
$ws.Cells.Item(13, 15).Value = 0.3055950511371977
$ws.Cells.Item(13, 16).Value = 0.3055950511371977
$ws.Cells.Item(13, 17).Value = 0.5611604862253333
$ws.Cells.Item(13, 18).Value = 5.050444376028
$ws.Cells.Item(13, 19).Value = 0.03369853271017832
$ws.Cells.Item(13, 20).Value = 0.03369853271017832

# Row 14
$ws.Cells.Item(14, 7).Value = 0.05795366666666666
$ws.Cells.Item(14, 8).Value = 0.173861
$ws.Cells.Item(14, 9).Value = 0.1102718535027888
$ws.Cells.Item(14, 10).Value = 0.1102718535027888
$ws.Cells.Item(14, 13).Value = 3.905093666666666
$ws.Cells.Item(14, 14).Value = 11.715281
$ws.Cells.Item(14, 15).Value = 0.1232456523180152
$ws.Cells.Item(14, 16).Value = 0.1232456523180152
$ws.Cells.Item(14, 17).Value = 0.2263144966601111
$ws.Cells.Item(14, 18).Value = 2.036830469941
$ws.Cells.Item(14, 19).Value = 0.01359052651726782
$ws.Cells.Item(14, 20).Value = 0.01359052651726782

# Row 15
$ws.Cells.Item(15, 7).Value = 0.05795366666666666
$ws.Cells.Item(15, 8).Value = 0.173861
$ws.Cells.Item(15, 9).Value = 0.1102718535027888
$ws.Cells.Item(15, 10).Value = 0.1102718535027888
$ws.Cells.Item(15, 13).Value = 6.285238333333333
$ws.Cells.Item(15, 14).Value = 18.855715
$ws.Cells.Item(15, 15).Value = 0.1983635642284282
$ws.Cells.Item(15, 16).Value = 0.1983635642284282
$ws.Cells.Item(15, 17).Value = 0.3642526072905555
$ws.Cells.Item(15, 18).Value = 3.278273465615
$ws.Cells.Item(15, 19).Value = 0.02187391789488827
$ws.Cells.Item(15, 20).Value = 0.02187391789488827

# Row 16
$ws.Cells.Item(16, 7).Value = 0.05795366666666666
$ws.Cells.Item(16, 8).Value = 0.173861
$ws.Cells.Item(16, 9).Value = 0.1102718535027888
$ws.Cells.Item(16, 10).Value = 0.1102718535027888
$ws.Cells.Item(16, 13).Value = 2.710232666666667
$ws.Cells.Item(16, 14).Value = 8.130698000000001
$ws.Cells.Item(16, 15).Value = 0.08553556494383548
$ws.Cells.Item(16, 16).Value = 0.08553556494383548
$ws.Cells.Item(16, 17).Value = 0.1570679205531111
$ws.Cells.Item(16, 18).Value = 1.413611284978
$ws.Cells.Item(16, 19).Value = 0.009432165286764903
$ws.Cells.Item(16, 20).Value = 0.009432165286764903
